$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TODO CMS")

# Row 18: client-side input checking implemented -> this item is no longer
# "nicht getestet" (untested), it is now "offen" (open) again.
$ws.Range("B18").Value2 = "offen"
$ws.Range("B15").Copy()
$ws.Range("B18").PasteSpecial(-4122) # xlPasteFormats

# New row 25: server now recognizes duplicate products -> add a new TODO
# item to test the client-side input checking, status "offen".
$ws.Range("A25").Value2 = "Eingabe testen"
$ws.Range("B25").Value2 = "offen"
$ws.Range("B15").Copy()
$ws.Range("B25").PasteSpecial(-4122) # xlPasteFormats

# Update the active selection to B18 (matches the saved workbook view).
$ws.Activate()
$ws.Range("B18").Select()
